# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2
# of the "zh-cn" and "de-de" worksheets to reflect the new report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 21:15:17"
$wsZhCn.Range("H2").Value = "2016-03-13 21:15:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 21:15:21"
$wsDeDe.Range("H2").Value = "2016-03-13 21:15:38"
